# Week3.pptx edit: fix "5: 30 AM" -> "5: 30 PM" typo on the "Timings for
# work" bullet (slide 9), and refresh the cached date placeholder text
# (5/20/22 -> 5/24/22) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

function Update-DateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $txt = $shp.TextFrame.TextRange.Text
            if ($txt -eq "5/20/22") {
                $shp.TextFrame.TextRange.Text = "5/24/22"
            }
        }
    }
}

# Slide master date placeholder.
Update-DateShape $p.SlideMaster

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li)
}

# Fix the AM/PM typo on slide 9, "Content Placeholder 2" shape, last
# paragraph ("Timings for work: Monday to Friday,  9:30 AM - 5: 30 AM").
$slide = $p.Slides.Item(9)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count
$para = $tr.Paragraphs($paraCount)
$oldText = $para.Text
$needle = "30 AM"
$startPos = $oldText.Length - $needle.Length + 1
$chars = $para.Characters($startPos, $needle.Length)
$chars.Text = "30 PM"
